$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "_id"
$ws.Range("B2").Value = 140

$ws.Range("C2").Select()
